$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Header text updates (new report week: 12/22/2025 - 12/28/2025, issue 52)
# ---------------------------------------------------------------------------
$a8 = $ws.Range("A8")
$a8.Value = $a8.Value().Replace("51", "52")

$c9 = $ws.Range("C9")
$c9v = $c9.Value()
$c9v = $c9v.Replace("12/15/2025", "12/22/2025")
$c9v = $c9v.Replace("12/21/2025", "12/28/2025")
$c9.Value = $c9v

# ---------------------------------------------------------------------------
# 2. New crime data for the week - update the stats table
# ---------------------------------------------------------------------------

# --- Row 15 (Rape) ---
$ws.Range("L15").Value = 0

# --- Row 16 (Robbery): prior-week comparison now has a zero base (N/A %) ---
$ws.Range("C15").Copy($ws.Range("D16"))
$ws.Range("E15").Copy($ws.Range("E16"))
$ws.Range("N16").Value = -94.607843137254

# --- Row 17 (Fel. Assault) ---
$ws.Range("I15").Copy($ws.Range("C17"))
$ws.Range("C17").Value = 1
$ws.Range("F17").Value = 3
$ws.Range("I17").Value = 22
$ws.Range("K17").Value = 83.333333333333
$ws.Range("L17").Value = 144.444444444444
$ws.Range("M17").Value = 266.666666666667
$ws.Range("N17").Value = -40.540540540540

# --- Row 19 (Gr. Larceny) ---
$ws.Range("L19").Value = -10.204081632653

# --- Row 21 (TOTAL, bold) ---
# Stash the original C21 (text, bold) formatting before it gets overwritten.
$ws.Range("C21").Copy($ws.Range("ZZ1"))

$ws.Range("C46").Copy($ws.Range("C21"))
$ws.Range("C21").Value = 1

$ws.Range("ZZ1").Copy($ws.Range("D21"))

$ws.Range("ZZ1").Copy($ws.Range("E21"))
$ws.Range("E21").Value = "***.*"

$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 3
$ws.Range("H21").Value = 66.666666666666
$ws.Range("I21").Value = 82
$ws.Range("J21").Value = 100
$ws.Range("K21").Value = -18
$ws.Range("L21").Value = -7.865168539325
$ws.Range("M21").Value = -18.811881188118
$ws.Range("N21").Value = -82.736842105263

$ws.Range("ZZ1").ClearContents()

# --- Row 24 (Petit Larceny) ---
$ws.Range("C15").Copy($ws.Range("C24"))
$ws.Range("F24").Value = 1
$ws.Range("H24").Value = 0
$ws.Range("L24").Value = -9.523809523809

# --- Row 26 (Misd. Assault) ---
$ws.Range("I15").Copy($ws.Range("D26"))
$ws.Range("D26").Value = 2
$ws.Range("K19").Copy($ws.Range("E26"))
$ws.Range("E26").Value = -100
$ws.Range("G26").Value = 3
$ws.Range("J26").Value = 38
$ws.Range("K26").Value = -31.578947368421

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("L27").Value = 66.666666666666

# --- Row 28 (Shooting Vic.) ---
$ws.Range("I15").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 1
$ws.Range("C15").Copy($ws.Range("D28"))
$ws.Range("E15").Copy($ws.Range("E28"))
$ws.Range("I15").Copy($ws.Range("F28"))
$ws.Range("F28").Value = 1
$ws.Range("H28").Value = -50
$ws.Range("I28").Value = 28
$ws.Range("K28").Value = 115.384615384615
$ws.Range("L28").Value = 55.555555555555
